# Append 4 new daily-attendance rows (72-75) to the Vacations_History sheet,
# following the exact same pattern as the existing rows (row 71 is the last
# populated row in the source workbook).
#
# Column layout (row 1 is the header row):
#   A = running day number (n-1 relative to row index)
#   B = date, stored as literal text (e.g. "2024-06-19"), NOT an Excel date
#   C = Arabic weekday name
#   D..G, I..AC = "موجود" (present)
#   H = "إمتداد" (extension)
#
# NOTE: Assigning a date-looking string straight to `.Value` makes Excel
# auto-convert it into a real date serial number, which is not what the
# source file does (it keeps plain text). To avoid that, the cell's
# NumberFormat is temporarily switched to Text ("@") before the value is
# written, and the original (General) formatting is restored immediately
# afterwards by re-pasting the formats copied from the template row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 72; Day = 71; Date = "2024-06-19"; Weekday = "الأربعاء" },
    @{ Row = 73; Day = 72; Date = "2024-06-20"; Weekday = "الخميس" },
    @{ Row = 74; Day = 73; Date = "2024-06-21"; Weekday = "الجمعة" },
    @{ Row = 75; Day = 74; Date = "2024-06-22"; Weekday = "السبت" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $src = $r - 1

    # Copy the whole template row's formatting (styles/borders/number
    # formats) onto the new row first.
    $ws.Range("A" + $src + ":AC" + $src).Copy()
    $ws.Range("A" + $r + ":AC" + $r).PasteSpecial(-4122)

    # Column A: sequential day number (plain integer).
    $ws.Cells.Item($r, 1).Value = $entry.Day

    # Column B: the date as literal text, not an auto-converted date value.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $entry.Date
    $ws.Range("A" + $src + ":B" + $src).Copy()
    $ws.Range("A" + $r + ":B" + $r).PasteSpecial(-4122)

    # Column C: Arabic weekday name.
    $ws.Cells.Item($r, 3).Value = $entry.Weekday

    # Columns D-G: present.
    $ws.Cells.Item($r, 4).Value = "موجود"
    $ws.Cells.Item($r, 5).Value = "موجود"
    $ws.Cells.Item($r, 6).Value = "موجود"
    $ws.Cells.Item($r, 7).Value = "موجود"

    # Column H: extension.
    $ws.Cells.Item($r, 8).Value = "إمتداد"

    # Columns I-AC: present.
    for ($col = 9; $col -le 29; $col++) {
        $ws.Cells.Item($r, $col).Value = "موجود"
    }
}
